$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = [double]"0.9967919210433758"
$ws.Range("E2").Value = [double]"0.9967919210433758"

$ws.Range("D3").Value = [double]"1.901561008950686E-11"
$ws.Range("E3").Value = [double]"1.901561008950686E-11"

$ws.Range("D4").Value = [double]"0.9900387634335711"
$ws.Range("E4").Value = [double]"0.9900387634335711"

$ws.Range("D5").Value = [double]"1.1062807477891E-06"
$ws.Range("E5").Value = [double]"1.1062807477891E-06"

$ws.Range("D6").Value = [double]"7.294668034478011E-31"
$ws.Range("E6").Value = [double]"7.294668034478011E-31"

$ws.Range("D7").Value = [double]"1.005831565672557E-17"
$ws.Range("E7").Value = 1

$ws.Range("D8").Value = [double]"0.9999999423688013"
$ws.Range("E8").Value = [double]"5.763119870216826E-08"

$ws.Range("D9").Value = [double]"0.9905236405392707"
$ws.Range("E9").Value = [double]"0.009476359460729333"

$ws.Range("D10").Value = [double]"5.550305138335896E-08"
$ws.Range("E10").Value = [double]"0.9999999444969486"

$ws.Range("D11").Value = [double]"0.001586649997082145"
$ws.Range("E11").Value = [double]"0.9984133500029179"
$ws.Range("F11").Value = [double]"7.26517391204834"
